# Books DB update: remove the "printit.txt" row and append newly added book entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 41 holds "printit.txt" -- delete the whole row so everything below shifts up.
$ws.Rows(41).Delete()

# New book cover images appended at the bottom of the list (rows 65-89).
$newBooks = @(
    "WhenGeniusFailed.jpg",
    "Lying.jpg",
    "TippingPoint.jpg",
    "CompetitiveAdvantage.jpg",
    "StrategicRiskTaking.jpg",
    "PyschopathWhisperer.jpg",
    "MalcomX.jpg",
    "PortableFinancialAnalyst.jpg",
    "FlashBoys.jpg",
    "Boomerang.jpg",
    "BigShort.jpg",
    "BlindWatchmaker.jpg",
    "Liar'sPoker.jpg",
    "WinningTheLoser'sGame.jpg",
    "MoneyChangesEverything.jpg",
    "Seveneves.jpg",
    "DeathByBlackHole.jpg",
    "StuffMatters.jpg",
    "ForTheLoveOfPhysics.jpg",
    "InvestingPsychology.jpg",
    "NarrativeAndNumbers.jpg",
    "PioneeringPortfolioManagement.jpg",
    "ASOIAF.jpg",
    "GunsGermsSteel.jpg",
    "HouseOfMorgan.jpg"
)

$startRow = 65
for ($i = 0; $i -lt $newBooks.Count; $i++) {
    $r = $startRow + $i
    $ws.Range("A$r").Value = $r
    $ws.Range("B$r").Value = "TitleHere"
    $ws.Range("C$r").Value = $newBooks[$i]
    $ws.Range("D$r").Value = "Brief description entered here"
    $ws.Range("E$r").Value = "LongDescription here"
}

# Update the sheet view to match where the author had scrolled to / selected.
$ws.Range("D72").Select()
